$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIEfIE")
$ws.Range("B2").Value = 1
